# RPA datasets push 2024-04-04
#
# A new IPO record ("미래" / 아이엠비디엑스) was added to the dataset as
# row 9, pushing every later row down by one; the former row 14
# (하나31호스팩) also now sorts after the former row 15 (에이피알/하나)
# entry. Rewrite rows 9-18 in place with the final, correct contents
# rather than juggling inserts, which keeps the shared-string table and
# cell styling byte-for-byte consistent with a plain re-pull of the
# source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final contents for every data row from row 9 (the newly inserted
# record) through the new last row (18).
$data = @(
    @(9,  "미래",  "2024-03-25", "아이엠비디엑스",   "미래",   "미래",      "2024-03-28", "2024-04-03", 32500, 2500000, 13000, 0, 100),
    @(10, "신한",  "2024-02-14", "에이피알",         "신한",   "신한, 하나", "2024-02-19", "2024-02-27", 75800, 379000,  250000, 0, 80),
    @(11, "유안타","2024-02-20", "유안타제15호스팩", "유안타", "유안타",    "2024-02-23", "2024-02-29", 13000, 6500000, 2000,   0, 100),
    @(12, "유진",  "2024-02-19", "유진스팩10호",     "유진",   "유진",      "2024-02-22", "2024-02-29", 8000,  4000000, 2000,   0, 100),
    @(13, "키움",  "2024-02-13", "코셈",             "키움",   "키움",      "2024-02-16", "2024-02-23", 9600,  600000,  16000,  0, 100),
    @(14, "하나",  "2024-03-18", "하나32호스팩",     "하나",   "하나",      "2024-03-21", "2024-03-27", 6000,  3000000, 2000,   0, 100),
    @(15, "하나",  "2024-02-14", "에이피알",         "신한",   "신한, 하나", "2024-02-19", "2024-02-27", 18950, 379000,  250000, 0, 20),
    @(16, "하나",  "2024-02-22", "하나31호스팩",     "하나",   "하나",      "2024-02-27", "2024-03-05", 10000, 5000000, 2000,   0, 100),
    @(17, "한국",  "2024-03-12", "삼현",             "한국",   "한국",      "2024-03-15", "2024-03-21", 60000, 2000000, 30000,  0, 100),
    @(18, "한화",  "2024-02-13", "이에이트",         "한화",   "한화",      "2024-02-16", "2024-02-23", 22600, 1130000, 20000,  0, 100)
)

# Columns that hold yyyy-mm-dd text in this sheet. Assigning a plain
# "2024-03-25"-style string through .Value auto-detects it as a date and
# reformats/serializes the cell, which the source workbook never does
# (every data cell here is plain text/number with the default style).
# Forcing a text number format first — then clearing the format back to
# the default "Normal" style once the literal text is safely stored —
# keeps the cell byte-identical to how the sheet already represents
# dates (shared-string text, no explicit style index).
$dateCols = @(2, 6, 7)

foreach ($row in $data) {
    $r = $row[0]
    for ($col = 1; $col -le 12; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $value = $row[$col]
        if ($dateCols -contains $col) {
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
